$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / extend the data table (A1:B10) ---------------------------
$data = @(
    @(2003, 453),
    @(2004, 421),
    @(2005, 360),
    @(2006, 497),
    @(2007, 231),
    @(2008, 460),
    @(2009, 422),
    @(2010, 288),
    @(2011, 448),
    @(2012, 287)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# --- Update the chart ----------------------------------------------------
$co = $ws.ChartObjects(1)
$chart = $co.Chart

# Chart style 4 -> 6
$chart.ChartStyle = 6

# Extend the series references from row 5 to row 10
$series = $chart.SeriesCollection(1)
$series.XValues = "=Sheet1!`$A`$1:`$A`$10"
$series.Values = "=Sheet1!`$B`$1:`$B`$10"

# Manual plot-area layout (factor mode, 85% width/height)
$plotArea = $chart.PlotArea
$plotArea.Width = 0.85
$plotArea.Height = 0.85
